$wb = $excel.ActiveWorkbook

# Loan RBI, Variable Instalments:
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / "heading" / "Outstanding" columns
# one place to the right, and leave that sheet active/selected with the
# cursor on K17 (mirrors the manual edit captured in the commit).

$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

$ws.Columns("N").Insert()

# Excel carries the left-neighbour column's width onto a freshly inserted
# column; column M is 11 characters wide, so column N should come out the
# same (minus the auto bestFit flag, since it is now blank).
$ws.Columns("N").ColumnWidth = 10.1666666666667

[void]$ws.Range("K17").Select()
